# update J9 part number and BOM
#
# Adds the J9 connector (a 4-pin header, same footprint as J4 "Conn_01x04")
# as a new line item on both the "BOM" and "DK Order" sheets, and restores
# the sheet/selection UI state recorded in the saved file.

$wb = $excel.ActiveWorkbook

$bom = $wb.Worksheets.Item("BOM")
$dk  = $wb.Worksheets.Item("DK Order")

# --- BOM sheet: new row 47 -------------------------------------------------
$bom.Range("A47").Value = 1
$bom.Range("B47").Value = "J9"
$bom.Range("C47").Value = "Conn_01x04"
$bom.Range("D47").Value = "S5596-ND"
$bom.Range("E47").Value = "NPTC041KFXC-RC"

# --- DK Order sheet: new row 47 --------------------------------------------
$dk.Range("A47").Value = 1
$dk.Range("B47").Value = "J9"
$dk.Range("C47").Value = "S5596-ND"

# --- UI state ---------------------------------------------------------------
# DK Order keeps a remembered selection even though it is no longer the
# active tab.
$dk.Range("C55").Select()

# BOM becomes the active/selected tab, with its own remembered selection.
$bom.Activate()
$bom.Range("E58").Select()
